$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F column (C_HCl) was a text "NaN" placeholder; replace with the actual
# measured value (10) for every data row.
$ws.Range("F2:F8").Value = 10

# Newly-read mass parameters from the master titration file: slope_NaOH
# (J) and intercept_NaOH (K) for each sample row.
$ws.Range("J2:J8").Value = 0.21
$ws.Range("K2:K8").Value = 3.5

# The trailing filename-list scratch rows (20-22) are no longer needed now
# that the values are read directly from the master titration file.
$ws.Rows("20:22").Delete()

# Update the active selection to match the saved view.
$ws.Range("I25").Select()
